$wb = $excel.ActiveWorkbook

# --- semantic_aspect_model_schema sheet ---
$schemaSheet = $wb.Worksheets.Item("semantic_aspect_model_schema")

# Widen column A (was narrow spacer width 2.4 -> 9.6 to fit the renamed header)
$schemaSheet.Columns.Item(1).ColumnWidth = 9.6

# Rename header "id" -> "dtwin_id" to avoid clashing with aspect model column names
$schemaSheet.Range("A1").Value = "dtwin_id"

# --- description sheet ---
$descriptionSheet = $wb.Worksheets.Item("description")
$descriptionSheet.Range("A5").Value = "dtwin_id"

# --- metadata sheet ---
$metadataSheet = $wb.Worksheets.Item("metadata")
$metadataSheet.Range("B2").Value = "41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$metadataSheet.Range("B3").Value = "https://github.com/dataspacesolutions/sldt-semantic-models/commit/41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$metadataSheet.Range("B4").Value = "2025-03-10 14:48:29+00:00"
